$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data for existing rows 2-5 (Q0..Q3) and new rows 6-11 (Q4..Q9)
# Columns: A=label, B..F=metrics, G=N
$labels = @("Q0", "Q1", "Q2", "Q3", "Q4", "Q5", "Q6", "Q7", "Q8", "Q9")

$data = @(
    @(0.06224550739541238, 0.4815100860261196, 0.5188979745415299, 0.7203457326461578, 0.7255814477924799, 46),
    @(0.009644723670651665, 0.4809544651665302, 0.5144163341041081, 0.7172282301360622, 0.725267178095959, 45),
    @(0.04605011742652749, 0.4545971144787497, 0.4787573890230339, 0.691922964659386, 0.6983704975779645, 44),
    @(0.02000496247511335, 0.4699992430775199, 0.4788062426551278, 0.6919582665559592, 0.6998547450934846, 43),
    @(0.05983179997382334, 0.4610733700736253, 0.4940437835170411, 0.702882482010358, 0.7088204810331585, 42),
    @(0.02666374005585367, 0.4762874528960522, 0.4947189357917796, 0.7033625919764142, 0.7115884908266601, 41),
    @(0.05161238290339588, 0.4630942636991037, 0.5133492645108798, 0.7164839597024345, 0.7237264155437525, 40),
    @(0.03337024079663585, 0.4821917770522179, 0.5171250100704028, 0.7191140452462341, 0.7277298162397244, 39),
    @(0.06108869437512898, 0.4568004358409793, 0.5137306737460197, 0.7167500776044741, 0.7237282623592945, 38),
    @(0.03255594249262348, 0.4650544728812123, 0.502363057293461, 0.708775745418437, 0.7177940205234036, 37)
)

# Copy the existing formatted label cell (A2) so new label cells (A6:A11) match
# the same style (border, bold, centered) without introducing new style entries.
$ws.Range("A2").Copy()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]

    if ($row -gt 5) {
        $ws.Range("A$row").PasteSpecial(-4122)
    }
    $ws.Cells.Item($row, 1).Value = $labels[$i]

    $ws.Cells.Item($row, 2).Value = $rowData[0]
    $ws.Cells.Item($row, 3).Value = $rowData[1]
    $ws.Cells.Item($row, 4).Value = $rowData[2]
    $ws.Cells.Item($row, 5).Value = $rowData[3]
    $ws.Cells.Item($row, 6).Value = $rowData[4]
    $ws.Cells.Item($row, 7).Value = $rowData[5]
}
